# "Doing what mosi wants" - trim the bullets sheet down to columns A:D and
# update a few of the remaining answers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused columns E:H entirely (dimension becomes A1:D7).
$ws.Range("E1:H7").EntireColumn.Delete()

# Row 3 - replace the remaining answer and clear the other two cells.
$ws.Range("B3").Value = "this is a very very very long question that cannot be answered easily by any person living on earth"
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()

# Row 4 - clear all three answer cells.
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()

# Row 5 - replace the remaining answer and clear the other two cells.
$ws.Range("B5").Value = "OK"
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()

# Row 6 - clear all three answer cells.
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()

# Row 7 - clear all three answer cells.
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
